$wb = $excel.ActiveWorkbook

# --- Sheet: Productdata (C2:C28, E2:E28) ---
$ws = $wb.Worksheets.Item("Productdata")

$cValues = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 5, 5, 5, 5, 5, 5, 5, 5)
$eValues = @(
    0.03401999999999999,
    0.025128,
    0.025254,
    0.024975,
    0.027999,
    0.0288,
    0.002880000000000001,
    0.004464000000000001,
    0.00432,
    0.00306,
    0.002988,
    0.004284,
    0.0009270000000000001,
    0.025218,
    0.0009630000000000001,
    0.002988,
    0.001143,
    0.032346,
    0.0007920000000000001,
    0.0131625,
    0.0132795,
    0.0131985,
    0.0135675,
    0.0146475,
    0.014562,
    0.0151155,
    0.015291
)

for ($i = 0; $i -lt $cValues.Length; $i++) {
    $row = $i + 2
    $ws.Range("C$row").Value = $cValues[$i]
    $ws.Range("E$row").Value = $eValues[$i]
}

# --- Sheet: Capacity (B2:B28) ---
$ws = $wb.Worksheets.Item("Capacity")

$bValues = @(240, 160, 40, 20, 60, 20, 320, 160, 160, 160, 400, 320, 20, 120, 20, 80, 80, 40, 40, 10, 10, 30, 40, 30, 20, 10, 40)

for ($i = 0; $i -lt $bValues.Length; $i++) {
    $row = $i + 2
    $ws.Range("B$row").Value = $bValues[$i]
}

# --- Sheet: ProcessingTime (diagonal cells) ---
$ws = $wb.Worksheets.Item("ProcessingTime")

$ws.Range("B2").Value = 3
$ws.Range("D4").Value = 2
$ws.Range("E5").Value = 1
$ws.Range("F6").Value = 3
$ws.Range("G7").Value = 1
$ws.Range("H8").Value = 4
$ws.Range("J10").Value = 2
$ws.Range("L12").Value = 5
$ws.Range("M13").Value = 4
$ws.Range("O15").Value = 3
$ws.Range("P16").Value = 1
$ws.Range("Q17").Value = 1
$ws.Range("S19").Value = 1
$ws.Range("T20").Value = 2
$ws.Range("V22").Value = 1
$ws.Range("W23").Value = 3
$ws.Range("X24").Value = 4
$ws.Range("Y25").Value = 3
$ws.Range("Z26").Value = 2
$ws.Range("AA27").Value = 1
$ws.Range("AB28").Value = 4
